$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" data column (H), mirroring column G's formatting,
# and fill in the values for each data row (4-25).
$values = @{
    4  = 2020
    5  = 42.2
    7  = 42.5
    8  = 42
    10 = 50.9
    11 = 36.9
    12 = 34.799999999999997
    14 = 30.7
    15 = 48.8
    17 = 61.1
    18 = 56.7
    19 = 41.6
    20 = 49
    21 = 43.5
    22 = 33.9
    23 = 34.6
    24 = 23.6
    25 = 35.9
}

# Copy column G's cell formatting into column H for every data row,
# including the blank rows (6, 9, 13, 16) which only need formatting.
foreach ($row in 4..25) {
    $ws.Range("G$row").Copy()
    $ws.Range("H$row").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Apply the values collected above.
foreach ($row in $values.Keys) {
    $ws.Range("H$row").Value = $values[$row]
}

# Rows 8 and 20 use the "0.0" number-format style (matching rows such as
# G10/G22/G23) instead of the plain style copied from column G.
$ws.Range("G10").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value = $values[8]

$ws.Range("G10").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H20").Value = $values[20]
$excel.CutCopyMode = $false

# Update the view: scroll back so column A is visible and move the
# active selection/cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B13").Select() | Out-Null
